$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the strain (marker) values for rows 47-49
$ws.Range("E47").Value = "TDY2188"
$ws.Range("E48").Value = "TDY2188"
$ws.Range("E49").Value = "TDY2188"

# Fill in marker_1 values for rows 47-52
$ws.Range("J47").Value = "NAT"
$ws.Range("J48").Value = "NAT"
$ws.Range("J49").Value = "NAT"
$ws.Range("J50").Value = "NAT"
$ws.Range("J51").Value = "NAT"
$ws.Range("J52").Value = "NAT"

# Update the active view window position / selection to reflect the new state
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L50").Select()
